$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 14; this shifts the existing rows 14-62 down to 15-63,
# preserving all of their data (Excel's native row-insert/shift behaviour).
$ws.Rows.Item(14).Insert()

# Populate the newly-inserted row 14 with a new weekly record. It duplicates the
# (now-shifted) row 15 data except for a new date (D=44883).
$ws.Cells.Item(14, 1).Value = 4
$ws.Cells.Item(14, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(14, 3).Value = "Los Lagos"
$ws.Cells.Item(14, 4).Value = 44883
$ws.Cells.Item(14, 5).Value = 10
$ws.Cells.Item(14, 6).Value = 300000000
$ws.Cells.Item(14, 7).Value = "Espárragos"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 600
$ws.Cells.Item(14, 11).Value = 1500
$ws.Cells.Item(14, 12).Value = 1700
$ws.Cells.Item(14, 13).Value = 1600
$ws.Cells.Item(14, 14).Value = "$/kilo"
$ws.Cells.Item(14, 15).Value = "Provincia de Linares"
$ws.Cells.Item(14, 16).Value = 1600
$ws.Cells.Item(14, 17).Value = 1
$ws.Cells.Item(14, 18).Value = "Hortaliza"

# Note: no explicit number-format fix-up needed for D14 — Excel's native
# row-insert already carries the date format (style index 2, numFmt
# "YYYY-MM-DD HH:MM:SS") down from the row above into the newly inserted row.
